$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 37; this shifts existing rows 37-72 down to 38-73.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record.
$ws.Cells.Item(37, 1).Value = 10
$ws.Cells.Item(37, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(37, 3).Value = "La Araucanía"
$ws.Cells.Item(37, 4).Value = 45040
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
$ws.Cells.Item(37, 5).Value = 9
$ws.Cells.Item(37, 6).Value = 100112042
$ws.Cells.Item(37, 7).Value = "Locoto"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 120
$ws.Cells.Item(37, 11).Value = 4400
$ws.Cells.Item(37, 12).Value = 4400
$ws.Cells.Item(37, 13).Value = 4400
$ws.Cells.Item(37, 14).Value = "$/kilo"
$ws.Cells.Item(37, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(37, 16).Value = 4400
$ws.Cells.Item(37, 17).Value = 1
$ws.Cells.Item(37, 18).Value = "Hortaliza"
